$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Claim Filing")
$ws.Columns.Item(3).AutoFit()
